$d = $word.ActiveDocument

# Find the paragraph "Shoots in a pattern. More HP." which ends the "Boss"
# section; the new "Still need" checklist goes right after it.
$anchor = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq "Shoots in a pattern. More HP.") {
        $anchor = $para
        break
    }
}
if ($anchor -eq $null) {
    throw "Could not locate anchor paragraph 'Shoots in a pattern. More HP.'"
}

function Add-PlainParagraph($afterPara, $text) {
    $afterPara.Range.InsertParagraphAfter()
    $newPara = $afterPara.Next()
    $newPara.Range.Text = $text
    return $newPara
}

$anchor = Add-PlainParagraph $anchor "Still need: "
$anchor = Add-PlainParagraph $anchor "Sound"
$anchor = Add-PlainParagraph $anchor "Animations"
$anchor = Add-PlainParagraph $anchor "Menus"
$anchor = Add-PlainParagraph $anchor "Gameplay loop"

# "Score" carries a lastRenderedPageBreak before its run text, matching the
# pagination mark Word stamped when the author's document was last saved.
$anchor.Range.InsertParagraphAfter()
$anchor = $anchor.Next()
$scoreRange = $anchor.Range
$scoreRange.Collapse(1)
$null = $scoreRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:lastRenderedPageBreak/><w:t>Score</w:t></w:r></w:p>")
$anchor = $d.Paragraphs.Item($anchor.Index)

$anchor = Add-PlainParagraph $anchor "Obstacles"
$anchor = Add-PlainParagraph $anchor "Player hp/lives"

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
